$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-76 down to 28-77
$ws.Rows(27).Insert()

# Populate the newly inserted row 27 with this week's price report
$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Vega Monumental Concepción"
$ws.Range("C27").Value = "Bíobío"
$ws.Range("D27").Value = 44708
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 100112012
$ws.Range("G27").Value = "Espinaca"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 5500
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = 5750
$ws.Range("N27").Value = "$/cuna 10 kilos"
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 575
$ws.Range("Q27").Value = 10
$ws.Range("R27").Value = "Hortaliza"
